$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant rows to append at the bottom of the table (rows 181-184)
$data = @(
    @("Asqarova Zarnigor Asqarovna", "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik", "AA8561113", "774", "Surxondaryo viloyati", "Denov tumani", "998900711166", "27-11-2024", "+998900711166"),
    @("Shermetova Intizor Hasan qizi", "Defektologiya (logopediya) 576 soatlik", "AB5889957", "775", "Xorazm viloyati", "Shovot tumani", "998913809033", "27-11-2024", "+998913809033"),
    @("Abdullayeva Oltinoy Choriyevna", "Defektologiya (logopediya) 576 soatlik", "AB3378808", "776", "Surxondaryo viloyati", "Uzun tumani", "998948626560", "27-11-2024", "+998948626560"),
    @("Yuldashov Jaloliddin Sharofiddin o'g'li", "Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik", "AB0767338", "777", "Xorazm viloyati", "Hazorasp tumani", "998973633103", "28-11-2024", "+998973633103")
)

$startRow = 181
$endRow = $startRow + $data.Count - 1

# Force the numeric-looking columns (Shartnoma raqam, Telefon raqam) to be
# stored as plain text so values such as "774" or "+998900711166" are not
# silently reinterpreted as numbers (which would also drop the leading "+").
$ws.Range("D$startRow`:D$endRow").NumberFormat = "@"
$ws.Range("G$startRow`:G$endRow").NumberFormat = "@"
$ws.Range("I$startRow`:I$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
